$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "lermacasemiro979"
$ws.Range("B4").Value = "Luz"
$ws.Range("C4").Value = "Macías"
$ws.Range("D4").Value = "chita77@compania.com"
$ws.Range("E4").Value = "#r(5dVgl)Cd2"
$ws.Range("F4").Value = "#r(5dVgl)Cd2"
$ws.Range("G4").Value = "Válido"
